$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the backlog item text (shared strings used in A4/D11 and A5/D12)
$ws.Range("A4").Value = "Upload sprints to github (Todo agente)"
$ws.Range("D11").Value = "Upload sprints to github (Todo agente)"

$ws.Range("A5").Value = "Take metrics (Todo agente)"
$ws.Range("D12").Value = "Take metrics (Todo agente)"

# Update the selected/active cell on the sheet (burndown chart selection change)
$null = $ws.Range("B5").Select()
